$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Renumber the sub-tasks under section 1.2 from 1.2.4/1.2.5/1.2.6 to 1.2.1/1.2.2/1.2.3
$ws.Range("A10").Value = "1.2.1"
$ws.Range("A11").Value = "1.2.2"
$ws.Range("A12").Value = "1.2.3"

# Remove the stray "Dependant on 2.1" comment in I16
$ws.Range("I16").ClearContents()

# Update progress value for row 18 (2.4) from 0.9 to 0.99
$ws.Range("E18").Value = 0.99

# Update the view: scroll to A25 and select B15
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B15").Select()
